$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 41.006619
$ws.Range("H2").Value = 123.019857
$ws.Range("I2").Value = 0.9349081063755518
$ws.Range("J2").Value = 0.9349081063755517
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 169.629438
$ws.Range("N2").Value = 508.888314
$ws.Range("O2").Value = 0.7428377317484701
$ws.Range("P2").Value = 0.7428377317484702
$ws.Range("Q2").Value = 6955.929735250122
$ws.Range("R2").Value = 62603.3676172511
$ws.Range("S2").Value = 0.6944850171332723
$ws.Range("T2").Value = 0.6944850171332723

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 41.006619
$ws.Range("H3").Value = 123.019857
$ws.Range("I3").Value = 0.9349081063755518
$ws.Range("J3").Value = 0.9349081063755517
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 0.9848756666666668
$ws.Range("N3").Value = 2.954627
$ws.Range("O3").Value = 0.004312947180081616
$ws.Range("P3").Value = 0.004312947180081616
$ws.Range("Q3").Value = 40.38642122537101
$ws.Range("R3").Value = 363.477791028339
$ws.Range("S3").Value = 0.004032209281027879
$ws.Range("T3").Value = 0.004032209281027878

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 41.006619
$ws.Range("H4").Value = 123.019857
$ws.Range("I4").Value = 0.9349081063755518
$ws.Range("J4").Value = 0.9349081063755517
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 54.620752
$ws.Range("N4").Value = 163.862256
$ws.Range("O4").Value = 0.2391940691454494
$ws.Range("P4").Value = 0.2391940691454494
$ws.Range("Q4").Value = 2239.812366757488
$ws.Range("R4").Value = 20158.31130081739
$ws.Range("S4").Value = 0.2236244742410349
$ws.Range("T4").Value = 0.2236244742410348

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 41.006619
$ws.Range("H5").Value = 123.019857
$ws.Range("I5").Value = 0.9349081063755518
$ws.Range("J5").Value = 0.9349081063755517
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 3.118221666666667
$ws.Range("N5").Value = 9.354665000000001
$ws.Range("O5").Value = 0.01365525192599884
$ws.Range("P5").Value = 0.01365525192599884
$ws.Range("Q5").Value = 127.867727842545
$ws.Range("R5").Value = 1150.809550582905
$ws.Range("S5").Value = 0.01276640572021668
$ws.Range("T5").Value = 0.01276640572021668

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 2.068862333333333
$ws.Range("H6").Value = 6.206586999999999
$ws.Range("I6").Value = 0.0471679015138598
$ws.Range("J6").Value = 0.04716790151385979
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 169.629438
$ws.Range("N6").Value = 508.888314
$ws.Range("O6").Value = 0.7428377317484701
$ws.Range("P6").Value = 0.7428377317484702
$ws.Range("Q6").Value = 350.939954902702
$ws.Range("R6").Value = 3158.459594124317
$ws.Range("S6").Value = 0.03503809697189084
$ws.Range("T6").Value = 0.03503809697189084

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 2.068862333333333
$ws.Range("H7").Value = 6.206586999999999
$ws.Range("I7").Value = 0.0471679015138598
$ws.Range("J7").Value = 0.04716790151385979
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 0.9848756666666668
$ws.Range("N7").Value = 2.954627
$ws.Range("O7").Value = 0.004312947180081616
$ws.Range("P7").Value = 0.004312947180081616
$ws.Range("Q7").Value = 2.037572169783223
$ws.Range("R7").Value = 18.338149528049
$ws.Range("S7").Value = 0.000203432667824569
$ws.Range("T7").Value = 0.0002034326678245689

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 2.068862333333333
$ws.Range("H8").Value = 6.206586999999999
$ws.Range("I8").Value = 0.0471679015138598
$ws.Range("J8").Value = 0.04716790151385979
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 54.620752
$ws.Range("N8").Value = 163.862256
$ws.Range("O8").Value = 0.2391940691454494
$ws.Range("P8").Value = 0.2391940691454494
$ws.Range("Q8").Value = 113.0028164311413
$ws.Range("R8").Value = 1017.025347880272
$ws.Range("S8").Value = 0.01128228229615193
$ws.Range("T8").Value = 0.01128228229615192

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 2.068862333333333
$ws.Range("H9").Value = 6.206586999999999
$ws.Range("I9").Value = 0.0471679015138598
$ws.Range("J9").Value = 0.04716790151385979
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 3.118221666666667
$ws.Range("N9").Value = 9.354665000000001
$ws.Range("O9").Value = 0.01365525192599884
$ws.Range("P9").Value = 0.01365525192599884
$ws.Range("Q9").Value = 6.451171353150555
$ws.Range("R9").Value = 58.06054217835499
$ws.Range("S9").Value = 0.0006440895779924577
$ws.Range("T9").Value = 0.0006440895779924575

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 0.7861760000000001
$ws.Range("H10").Value = 2.358528
$ws.Range("I10").Value = 0.01792399211058844
$ws.Range("J10").Value = 0.01792399211058843
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 169.629438
$ws.Range("N10").Value = 508.888314
$ws.Range("O10").Value = 0.7428377317484701
$ws.Range("P10").Value = 0.7428377317484702
$ws.Range("Q10").Value = 133.358593049088
$ws.Range("R10").Value = 1200.227337441792
$ws.Range("S10").Value = 0.01331461764330699
$ws.Range("T10").Value = 0.01331461764330699

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 0.7861760000000001
$ws.Range("H11").Value = 2.358528
$ws.Range("I11").Value = 0.01792399211058844
$ws.Range("J11").Value = 0.01792399211058843
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 0.9848756666666668
$ws.Range("N11").Value = 2.954627
$ws.Range("O11").Value = 0.004312947180081616
$ws.Range("P11").Value = 0.004312947180081616
$ws.Range("Q11").Value = 0.7742856121173336
$ws.Range("R11").Value = 6.968570509056002
$ws.Range("S11").Value = 0.00007730523122916752
$ws.Range("T11").Value = 0.00007730523122916751

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 0.7861760000000001
$ws.Range("H12").Value = 2.358528
$ws.Range("I12").Value = 0.01792399211058844
$ws.Range("J12").Value = 0.01792399211058843
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 54.620752
$ws.Range("N12").Value = 163.862256
$ws.Range("O12").Value = 0.2391940691454494
$ws.Range("P12").Value = 0.2391940691454494
$ws.Range("Q12").Value = 42.94152432435201
$ws.Range("R12").Value = 386.473718919168
$ws.Range("S12").Value = 0.00428731260826258
$ws.Range("T12").Value = 0.004287312608262579

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 0.7861760000000001
$ws.Range("H13").Value = 2.358528
$ws.Range("I13").Value = 0.01792399211058844
$ws.Range("J13").Value = 0.01792399211058843
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 3.118221666666667
$ws.Range("N13").Value = 9.354665000000001
$ws.Range("O13").Value = 0.01365525192599884
$ws.Range("P13").Value = 0.01365525192599884
$ws.Range("Q13").Value = 2.451471037013334
$ws.Range("R13").Value = 22.06323933312
$ws.Range("S13").Value = 0.0002447566277897008
$ws.Range("T13").Value = 0.0002447566277897008
